$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target data for the "Frame0" table (A1:C22); row 1 is the header
# (nazev | zkratka | seminariciUcitIdno) and stays untouched.
$rows = @(
    @("Podnikové informační systémy", "EIS", 14),
    @("Podnikové informační systémy", "KEIS", 14),
    @("Fyzikální praktikum C", "K521", 302),
    @("Počítačové modelování I", "K107", 612),
    @("Podnikové informační systémy", "EIS", 1609),
    @("Podnikové informační systémy", "KEIS", 1609),
    @("Podnikové informační systémy", "EIS", 3457),
    @("Podnikové informační systémy", "KEIS", 3457),
    @("Podnikové informační systémy", "EIS", 3606),
    @("Podnikové informační systémy", "KEIS", 3606),
    @("Sociální sítě", "SON", 4190),
    @("Praktické aplikace hardwaru", "AHW", 4746),
    @("Základy autonomní robotiky", "0182", 4746),
    @("Podnikové informační systémy", "EIS", 4991),
    @("Podnikové informační systémy", "KEIS", 4991),
    @("Reflektivní seminář pedagogické praxe", "KSPP", 8021),
    @("Softwarové inženýrství", "KSWI", 8093),
    @("Softwarové inženýrství", "SWI", 8093),
    @("Odborná prezentace", "KOPRE", 8514),
    @("Odborná prezentace", "OPRE", 8514),
    @("Introduction to MATLAB", "ITM", 8514)
)

# Make sure the two brand-new rows (21 and 22) inherit the same formatting
# (styles) as the existing data rows before filling them in.
$ws.Range("A20:C20").Copy() | Out-Null
$ws.Range("A21:C22").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# "0182" (row 14, column zkratka) is a course code that looks numeric. Mark
# that cell as Text first so Excel keeps the leading zero instead of
# silently turning it into the number 182.
$ws.Range("B14").NumberFormat = "@"

$rowIndex = 2
foreach ($r in $rows) {
    $ws.Cells.Item($rowIndex, 1).Value = $r[0]
    $ws.Cells.Item($rowIndex, 2).Value = $r[1]
    $ws.Cells.Item($rowIndex, 3).Value = $r[2]
    $rowIndex = $rowIndex + 1
}

# Grow the table ("Frame0") so it covers the two extra rows.
$tbl = $ws.ListObjects.Item(1)
$newRange = $ws.Range("A1:C22")
$tbl.Resize($newRange)

Write-Host "edit complete"
